$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate the last data row (16) down to a new row 17, preserving its
# formatting/styles, then overwrite the cell values with the new batch_016
# test case content.
$ws.Rows("16").Copy()
$ws.Rows("17").Insert(-4121)

# Remove the (empty) H17 cell that got inherited from H16; the new row
# should not have a value/style there.
$ws.Range("H17").Clear()

$ws.Range("A17").Value = "batch_016"
$ws.Range("B17").Value = "y"
$ws.Range("C17").Value = "批量操作语句16执行"
$ws.Range("D17").Value = "batchsql"
$ws.Range("E17").Value = "SingleTable"
$ws.Range("G17").Value = "M"
$ws.Range("I17").Value = "batch_sql_016"
$ws.Range("K17").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/batchsql/expectedresult/batch_016.csv"
$ws.Range("J17").Value = "select m.name as n from M as m order by m.name "
$ws.Range("N17").Value = "csv_containsAll"

# Nudge the number format so the style for G17 settles on the same style
# index used elsewhere in the workbook for this identical format.
$ws.Range("G17").NumberFormat = "@"

# Update the tracked selection, matching the workbook's last saved state.
$ws.Range("J22").Select()
